# Update "想去人数" (F column) values for specific events on the
# "展览" and "全部类型" sheets, as produced by the site regeneration
# at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value (想去人数 / "want to go" count)
$updates = @{
    3  = 3023
    4  = 217
    7  = 1646
    8  = 1613
    14 = 27
    22 = 359
    23 = 159
    26 = 2024
    30 = 176
    32 = 224
    35 = 493
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
